# Nits to variables 2 lesson/worksheet
$d = $word.ActiveDocument

function ReplaceText($findText, $replaceText) {
    $range = $d.Content
    $ok = $range.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)
    if (-not $ok) {
        Write-Output "NOT FOUND: $findText"
    }
    return $ok
}

# 1. "Explain in 1-2 sentences..." -- must -> should
ReplaceText "must be a space after the “is” in “Answer is: ”." "should be a space after the “is” in “Answer is: ”."

# 2. "...and then prints (says) “Hi NAME!”." -> "...and then says “Hi NAME!”."
ReplaceText "prints (says)" "says"

# 3. Rectangle paragraph restructuring
ReplaceText "of a rectangle and then prints the area of the rectangle. Then, if the area is smaller than 10, the program prints “" "of a rectangle. Then, it says the area of the rectangle. Finally, if the area is smaller than 10, the program says “"

# 4. "DO NOT draw a rectangle, this is just text." -> "DO NOT draw a rectangle; this is just text."
ReplaceText "draw a rectangle, this is just text." "draw a rectangle; this is just text."

# 5. circle paragraph -- prints out -> says
ReplaceText "radius of a circle and prints out the area of the circle" "radius of a circle and says the area of the circle"

# 6. "Then, if the area is bigger than 100, the program prints “" -> "...program says “"
ReplaceText "Then, if the area is bigger than 100, the program prints “" "Then, if the area is bigger than 100, the program says “"

# 7. "...and then prints the larger one." -> "...and then says the larger one."
ReplaceText "and then prints the larger one." "and then says the larger one."

# 8 & 9. Temperature paragraph restructuring
ReplaceText "asks the user for a temperature in Fahrenheit and converts it to Celsius. If the temperature in Celsius is greater than 100, the program prints “boiling”. If the number is less than 0, the program prints “freezing”." "asks the user for a temp in Fahrenheit and converts it to Celsius. Then, it says the temp in Celsius. Finally, if the temp in Celsius is greater than 100, the program says “boiling”. If the number is less than 0, the program says “freezing”."

# 10. "...and then prints “blast off” (hint: use a loop)." -> "...and then says “blast off” ..."
ReplaceText "and then prints “blast off”" "and then says “blast off”"

# 11. "Write a program the prints the largest of three numbers..." -> "...the says the largest..."
ReplaceText "Write a program the prints the largest of three numbers" "Write a program the says the largest of three numbers"

# 12. "Create a program that asks the user for 5 numbers then prints the sum..." -> "...for 8 numbers then says the sum..."
ReplaceText "asks the user for 5 numbers then prints the sum of all the numbers" "asks the user for 8 numbers then says the sum of all the numbers"

# 13. "Write a program that prints the number of spaces..." -> "...that says the number of spaces..."
ReplaceText "Write a program that prints the number of spaces in a sentence given by a user." "Write a program that says the number of spaces in a sentence given by a user."

# Now set the underline formatting on the "says" that was introduced in the rectangle paragraph
# (the one right after "Then, it " and before " the area of the rectangle.")
$sentence = $d.Content
$foundSentence = $sentence.Find.Execute("Then, it says the area of the rectangle.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($foundSentence) {
    $scoped = $d.Range($sentence.Start, $sentence.End)
    $foundWord = $scoped.Find.Execute("says", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($foundWord) {
        Write-Output "underline target text: [$($scoped.Text)]"
        $scoped.Font.Underline = 1
    } else {
        Write-Output "says not found inside sentence range"
    }
} else {
    Write-Output "sentence not found for underline step"
}
